$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1133093333333333
$ws.Range("H2").Value = 0.339928
$ws.Range("I2").Value = 0.02456654176752224
$ws.Range("J2").Value = 0.02456654176752224
$ws.Range("M2").Value = 133.3951123333333
$ws.Range("N2").Value = 400.185337
$ws.Range("O2").Value = 0.8984588679103155
$ws.Range("P2").Value = 0.8984588679103156
$ws.Range("Q2").Value = 15.11491124841511
$ws.Range("R2").Value = 136.034201235736
$ws.Range("S2").Value = 0.02207202730491951
$ws.Range("T2").Value = 0.02207202730491951

$ws.Range("G3").Value = 0.1133093333333333
$ws.Range("H3").Value = 0.339928
$ws.Range("I3").Value = 0.02456654176752224
$ws.Range("J3").Value = 0.02456654176752224
$ws.Range("M3").Value = 2.340788333333334
$ws.Range("N3").Value = 7.022365000000001
$ws.Range("O3").Value = 0.01576596023045448
$ws.Range("P3").Value = 0.01576596023045448
$ws.Range("Q3").Value = 0.2652331655244445
$ws.Range("R3").Value = 2.38709848972
$ws.Range("S3").Value = 0.0003873151205065544
$ws.Range("T3").Value = 0.0003873151205065546

$ws.Range("G4").Value = 0.1133093333333333
$ws.Range("H4").Value = 0.339928
$ws.Range("I4").Value = 0.02456654176752224
$ws.Range("J4").Value = 0.02456654176752224
$ws.Range("M4").Value = 12.735128
$ws.Range("N4").Value = 38.205384
$ws.Range("O4").Value = 0.08577517185923002
$ws.Range("P4").Value = 0.08577517185923003
$ws.Range("Q4").Value = 1.443008863594667
$ws.Range("R4").Value = 12.987079772352
$ws.Range("S4").Value = 0.002107199342096172
$ws.Range("T4").Value = 0.002107199342096173

$ws.Range("I5").Value = 0.8380577451911468
$ws.Range("J5").Value = 0.8380577451911468
$ws.Range("M5").Value = 133.3951123333333
$ws.Range("N5").Value = 400.185337
$ws.Range("O5").Value = 0.8984588679103155
$ws.Range("P5").Value = 0.8984588679103156
$ws.Range("Q5").Value = 515.6268456294275
$ws.Range("R5").Value = 4640.641610664847
$ws.Range("S5").Value = 0.7529604129879094
$ws.Range("T5").Value = 0.7529604129879095

$ws.Range("I6").Value = 0.8380577451911468
$ws.Range("J6").Value = 0.8380577451911468
$ws.Range("M6").Value = 2.340788333333334
$ws.Range("N6").Value = 7.022365000000001
$ws.Range("O6").Value = 0.01576596023045448
$ws.Range("P6").Value = 0.01576596023045448
$ws.Range("Q6").Value = 9.048107411812779
$ws.Range("S6").Value = 0.01321278508150797
$ws.Range("T6").Value = 0.01321278508150798

$ws.Range("I7").Value = 0.8380577451911468
$ws.Range("J7").Value = 0.8380577451911468
$ws.Range("M7").Value = 12.735128
$ws.Range("N7").Value = 38.205384
$ws.Range("O7").Value = 0.08577517185923002
$ws.Range("P7").Value = 0.08577517185923003
$ws.Range("Q7").Value = 49.22649536752267
$ws.Range("R7").Value = 443.038458307704
$ws.Range("S7").Value = 0.07188454712172943
$ws.Range("T7").Value = 0.07188454712172943

$ws.Range("G8").Value = 0.6336240000000001
$ws.Range("H8").Value = 1.900872
$ws.Range("I8").Value = 0.1373757130413309
$ws.Range("J8").Value = 0.1373757130413309
$ws.Range("M8").Value = 133.3951123333333
$ws.Range("N8").Value = 400.185337
$ws.Range("O8").Value = 0.8984588679103155
$ws.Range("P8").Value = 0.8984588679103156
$ws.Range("Q8").Value = 84.52234465709601
$ws.Range("R8").Value = 760.701101913864
$ws.Range("S8").Value = 0.1234264276174865
$ws.Range("T8").Value = 0.1234264276174866

$ws.Range("G9").Value = 0.6336240000000001
$ws.Range("H9").Value = 1.900872
$ws.Range("I9").Value = 0.1373757130413309
$ws.Range("J9").Value = 0.1373757130413309
$ws.Range("M9").Value = 2.340788333333334
$ws.Range("N9").Value = 7.022365000000001
$ws.Range("O9").Value = 0.01576596023045448
$ws.Range("P9").Value = 0.01576596023045448
$ws.Range("Q9").Value = 1.48317966692
$ws.Range("R9").Value = 13.34861700228
$ws.Range("S9").Value = 0.00216586002843995
$ws.Range("T9").Value = 0.00216586002843995

$ws.Range("G10").Value = 0.6336240000000001
$ws.Range("H10").Value = 1.900872
$ws.Range("I10").Value = 0.1373757130413309
$ws.Range("J10").Value = 0.1373757130413309
$ws.Range("M10").Value = 12.735128
$ws.Range("N10").Value = 38.205384
$ws.Range("O10").Value = 0.08577517185923002
$ws.Range("P10").Value = 0.08577517185923003
$ws.Range("Q10").Value = 8.069282743872002
$ws.Range("R10").Value = 72.62354469484801
$ws.Range("S10").Value = 0.01178342539540443
$ws.Range("T10").Value = 0.01178342539540443
